$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.968.50"
$ws.Range("E2").Value = '  +1.91%  '

$ws.Range("D3").Value = "'1.909.96"
$ws.Range("E3").Value = '  +2.53%  '

$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("E5").Value = '  +1.54%  '

$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("D7").Value = "'0.4815"
$ws.Range("E7").Value = '  +0.88%  '

$ws.Range("D8").Value = "'0.3808"
$ws.Range("E8").Value = '  +1.17%  '

$ws.Range("D9").Value = "'0.07350"
$ws.Range("E9").Value = '  +0.33%  '

$ws.Range("D10").Value = "'0.9322"
$ws.Range("E10").Value = '  -0.33%  '

$ws.Range("D11").Value = "'20.77"
$ws.Range("E11").Value = '  +0.27%  '

$ws.Range("E12").Value = '  -0.25%  '

$ws.Range("D13").Value = "'1.911.18"
$ws.Range("E13").Value = '  +2.48%  '

$ws.Range("D14").Value = "'5.494"
$ws.Range("E14").Value = '  +0.99%  '

$ws.Range("D15").Value = "'6.630"
$ws.Range("E15").Value = '  +1.12%  '

$ws.Range("D16").Value = "'91.87"
$ws.Range("E16").Value = '  +1.57%  '

$ws.Range("D17").Value = "'1.009"
$ws.Range("E17").Value = '  -0.27%  '

$ws.Range("D18").Value = "'0.000008844"
$ws.Range("E18").Value = '  -0.49%  '

$ws.Range("E19").Value = '  -0.41%  '

$ws.Range("D20").Value = "'28.011.03"
$ws.Range("E20").Value = '  +1.98%  '

$ws.Range("D21").Value = "'14.76"
$ws.Range("E21").Value = '  +0.30%  '

$ws.Range("D22").Value = "'5.165"
$ws.Range("E22").Value = '  +1.03%  '

$ws.Range("D23").Value = "'2.157.89"
$ws.Range("E23").Value = '  +1.73%  '

$ws.Range("D24").Value = "'10.90"
$ws.Range("E24").Value = '  +1.76%  '

$ws.Range("D25").Value = "'156.13"
$ws.Range("E25").Value = '  +0.40%  '

$ws.Range("D26").Value = "'1.913"
$ws.Range("E26").Value = '  -1.24%  '

$ws.Range("D27").Value = "'18.48"
$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("D28").Value = "'2.117"
$ws.Range("E28").Value = '  +4.43%  '

$ws.Range("D29").Value = "'116.77"
$ws.Range("E29").Value = '  +0.99%  '

$ws.Range("D30").Value = "'4.949"
$ws.Range("E30").Value = '  -0.03%  '

$ws.Range("D31").Value = "'0.08946"
$ws.Range("E31").Value = '  +0.64%  '

$ws.Range("E32").Value = '  -0.59%  '

$ws.Range("D33").Value = "'1.257"
$ws.Range("E33").Value = '  +3.59%  '

$ws.Range("D34").Value = "'0.7756"
$ws.Range("E34").Value = '  +2.24%  '

$ws.Range("D35").Value = "'4.667"
$ws.Range("E35").Value = '  +1.31%  '

$ws.Range("D36").Value = "'2.606"
$ws.Range("E36").Value = '  -4.55%  '

$ws.Range("D37").Value = "'0.02051"
$ws.Range("E37").Value = '  -0.17%  '

$ws.Range("D38").Value = "'1.108"
$ws.Range("E38").Value = '  -0.88%  '

$ws.Range("D39").Value = "'0.05305"
$ws.Range("E39").Value = '  +0.34%  '

$ws.Range("D40").Value = "'0.5501"
$ws.Range("E40").Value = '  -3.11%  '

$ws.Range("D41").Value = "'2.998"
$ws.Range("E41").Value = '  +0.39%  '

$ws.Range("D42").Value = "'7.016"
$ws.Range("E42").Value = '  -0.87%  '

$ws.Range("D43").Value = "'0.1525"
$ws.Range("E43").Value = '  -0.29%  '

$ws.Range("D44").Value = "'8.459"
$ws.Range("E44").Value = '  -2.66%  '

$ws.Range("D45").Value = "'10.65"
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").Value = "'108.51"
$ws.Range("E46").Value = '  +5.20%  '

$ws.Range("D47").Value = "'0.4812"
$ws.Range("E47").Value = '  -2.03%  '

$ws.Range("E48").Value = '  -0.53%  '

$ws.Range("D49").Value = "'1.642"
$ws.Range("E49").Value = '  -1.22%  '

$ws.Range("D50").Value = "'68.02"
$ws.Range("E50").Value = '  +0.96%  '

$ws.Range("D51").Value = "'0.06084"
$ws.Range("E51").Value = '  +0.12%  '

Write-Output "Applied cryptos update"
